# 970724 - +2232 - vbmellat part sell
#
# Update the "part sell" inputs on Sheet1/Sheet3 (the sale date moves one
# day later, the sold quantities/prices for the three lots change, and
# Sheet1's manual "C11" override becomes a plain number instead of a
# formula). Every other changed cell in the diff is a formula result that
# recalculates automatically once these raw inputs change.

$wb  = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# --- Sheet1: reference date moves from 35634 (1997-07-23) to 35635 (1997-07-24) ---
$ws1.Range("F7").Value = 35635

# Sheet1!C11 used to be the formula "=2400000"; it is now a plain value.
$ws1.Range("C11").Value = 17400000

# Sheet1!E11/F11 quantity+price inputs
$ws1.Range("E11").Value = 630
$ws1.Range("F11").Value = 850

# --- Sheet2: sold quantity for lot "50" changes ---
$ws2.Range("C15").Value = 37119

# --- Sheet3: sale date for all three lots moves one day later, and the
#     sold quantities change for each lot ---
$ws3.Range("B14").Value = 35635
$ws3.Range("D14").Value = 9895

$ws3.Range("B15").Value = 35635
$ws3.Range("D15").Value = 2450

$ws3.Range("B16").Value = 35635
$ws3.Range("D16").Value = 14499

# --- Selection / active sheet bookkeeping to match the saved UI state ---
$ws1.Activate()
$ws1.Range("C11").Select()

$ws2.Activate()
$ws2.Range("C16").Select()

$ws3.Activate()
$ws3.Range("D15").Select()
